# Add two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# existing header style used by the other headers in row 1, and fill in
# the data rows (2-38) with their corresponding values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - use same style as existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-38
$data = @(
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(1, 1),
    @(7, 7),
    @(7, 7),
    @(8, 9),
    @(6, 7),
    @(7, 8),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(5, 5),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(5, 6),
    @(3, 3),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(4, 4),
    @(5, 5),
    @(6, 7),
    @(7, 7),
    @(5, 5),
    @(9, 9),
    @(7, 7),
    @(5, 5),
    @(5, 5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
